# Update "countries & provincias Spain" COVID data table.
#
# Three provinces got updated case counts (Zaragoza, Teruel, Huesca), and
# the footer timestamp advanced from 14:16 to 14:46. Because the data
# table (A4:E61) is kept sorted descending by "Casos totales" (column B),
# bumping those three values changes their rank and re-shuffles the rows
# around them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp (A1).
$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 14:46"

# Zaragoza: 224/0/210/14 -> 278/0/264/14 (currently row 19).
$ws.Range("B19").Value = 278
$ws.Range("D19").Value = 264
$ws.Range("E19").Value = 14

# Teruel: 27/0/26/1 -> 40/0/38/2 (currently row 50).
$ws.Range("B50").Value = 40
$ws.Range("D50").Value = 38
$ws.Range("E50").Value = 2

# Huesca: 24/0/24/0 -> 34/0/34/0 (currently row 53).
$ws.Range("B53").Value = 34
$ws.Range("D53").Value = 34

# Re-sort the data range by "Casos totales" descending, same as the
# original table ordering, so the updated provinces land on their new rank.
$dataRange = $ws.Range("A4:E61")
$sortKey = $ws.Range("B4:B61")
$dataRange.Sort($sortKey, 2)

# La Palma and Arroyo de la Luz are tied at 7 total cases; the source
# listing swapped their relative order on this update even though the
# numbers match, so fix that up explicitly post-sort.
$first = $ws.Range("A56").Text
$second = $ws.Range("A57").Text
$ws.Range("A56").Value = $second
$ws.Range("A57").Value = $first
